# Ajuste de login page
# Update the timestamp embedded in the test-user e-mail addresses from
# "20251109_005042" to "20251109_011412" wherever it appears in the workbook.

$wb = $excel.ActiveWorkbook

$oldStamp = "20251109_005042"
$newStamp = "20251109_011412"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
            $cell.Value = $val.Replace($oldStamp, $newStamp)
        }
    }
}
